$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "I made a change!"

$ws.Range("J10").Select()
